$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A62: convert stored text "1388412" to a real number ---
$ws.Cells.Item(62, 1).Value = 1388412

# --- Append new row 63 (matches the source-data pattern: every cell is text) ---
$row = 63

# Job ID (kept as text, per the existing sheet convention for every column but A62)
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "1388543"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 'Graphics Designer / Digital Marketer'
$ws.Cells.Item($row, 3).Value = 'NGICON'
$ws.Cells.Item($row, 4).Value = 'Jul 23, 2025'
$ws.Cells.Item($row, 5).Value = 'Aug 2, 2025'

# Vacancies - numeric-looking text, force text storage like the rest of the column
$ws.Cells.Item($row, 6).NumberFormat = "@"
$ws.Cells.Item($row, 6).Value = '3'
$ws.Cells.Item($row, 6).Style = "Normal"

$ws.Cells.Item($row, 7).Value = 'Full Time'
$ws.Cells.Item($row, 8).Value = 'Work at office'
$ws.Cells.Item($row, 9).Value = 'Bachelor/Honors'
$ws.Cells.Item($row, 10).Value = '2 to 6 years'
$ws.Cells.Item($row, 11).Value = 'Age 18 to 50 years; Only Female; Skills Required: Creative The applicants should have experience in the following area(s): Brand Promotion, Campaign Planning, Content writingThe applicants should have experience in the following business area(s):E-commerce The applicants should have experience in the following area(s): Digital, Social Media Operation, Ability to Work Under Pressure, Advertisement, Branding, Campaign Planning, Marketing, Promotion Experience in E-Commerce Organization.Expertise on Social Media platforms marketing, Facebook, Instagram, LinkedIn, Google Ad-Sense, manage digital marketing campaigns, and YouTube marketing.Must be smart, good looking, have beautiful voice, expert on live presentation and willing to work on social media.Experienced candidates in live presentation in Facebook, modelling and sales will be preferred.Responsibilities & ContextDigital Marketing Executive - Live Presenter (Facebook, YouTube, Video Editing & Product Photography)Manage Facebook page (Example: Comments reply, Checking inbox) Manage large amounts of incoming calls, message and customer query.Host Live Event.Product Launches presentation.Must have idea to received customer order with name address and phone number.Can write down all sales and delivery data on record book.Plan and execute all the digital marketing campaign, including a good Graphics Design Knowledge, Photoshop & Illustrator, create necessary marketing database, email and social media.Facebook pixel, Facebook boosting, email marketing, SMS campaign. Excellent understanding of Facebook, Alibaba, LinkedIn, Twitter, Pinterest, Youtube, Instagram and Snapchat as marketing platforms.; Skills Required: Creative The applicants should have experience in the following area(s): Brand Promotion, Campaign Planning, Content writing; The applicants should have experience in the following business area(s):E-commerce The applicants should have experience in the following area(s): Digital, Social Media Operation, Ability to Work Under Pressure, Advertisement, Branding, Campaign Planning, Marketing, Promotion Experience in E-Commerce Organization.; Expertise on Social Media platforms marketing, Facebook, Instagram, LinkedIn, Google Ad-Sense, manage digital marketing campaigns, and YouTube marketing.; Must be smart, good looking, have beautiful voice, expert on live presentation and willing to work on social media.; Experienced candidates in live presentation in Facebook, modelling and sales will be preferred.; Responsibilities & Context; Digital Marketing Executive - Live Presenter (Facebook, YouTube, Video Editing & Product Photography); Manage Facebook page (Example: Comments reply, Checking inbox) Manage large amounts of incoming calls, message and customer query.; Host Live Event.; Product Launches presentation.; Must have idea to received customer order with name address and phone number.; Can write down all sales and delivery data on record book.; Plan and execute all the digital marketing campaign, including a good Graphics Design Knowledge, Photoshop & Illustrator, create necessary marketing database, email and social media.; Facebook pixel, Facebook boosting, email marketing, SMS campaign. Excellent understanding of Facebook, Alibaba, LinkedIn, Twitter, Pinterest, Youtube, Instagram and Snapchat as marketing platforms.'
$ws.Cells.Item($row, 12).Value = 'Adobe illustrator; Adobe Photoshop; Adobe Premiere Pro; DaVinchi Resolve; Digital Content Development; Digital Marketing; Digital Marketing (Social Media Marketing); Digital marketing expert; SEO and Digital Marketer; Video Editing.'
$ws.Cells.Item($row, 13).Value = 'Make designs social media posts and ads.; Make designs social media posts and ads.; Create materials for events and promotions Video Editing for Facebook and YouTube.; Create materials for events and promotions Video Editing for Facebook and YouTube.; Design for required projects (such as brochure, leaflet, banner, logo, poster etc.); Design for required projects (such as brochure, leaflet, banner, logo, poster etc.); Stay current with design trends, tools, and software to ensure our visual content remains fresh and innovative.; Stay current with design trends, tools, and software to ensure our visual content remains fresh and innovative.; Handle camera on some occasions and capture footage as and when required for the video content.; Handle camera on some occasions and capture footage as and when required for the video content.; Create social media posts, Google Ads banners, website elements, YouTube videos, animations, presentations, reels, shorts and more.; Create social media posts, Google Ads banners, website elements, YouTube videos, animations, presentations, reels, shorts and more.'
$ws.Cells.Item($row, 14).Value = 'Dhaka (DOHS Mirpur)'
$ws.Cells.Item($row, 15).Value = 'Tk. 16000 - 32000 (Monthly)'
$ws.Cells.Item($row, 16).Value = 'House-343(2nd floor), Avenue-3, Road-5, Mirpur DOHS'
$ws.Cells.Item($row, 17).Value = ""
$ws.Cells.Item($row, 18).Value = ""
$ws.Cells.Item($row, 19).Value = 'https://jobs.bdjobs.com/jobdetails.asp?id=1388543'
